$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# Step 1: insert a new "2022-Q1" sheet right before the "总计" sheet,
# mirroring the per-fund layout used by "2021-Q2".
# ------------------------------------------------------------------
$src        = $wb.Worksheets.Item("2021-Q2")
$zongjiRef  = $wb.Worksheets.Item("总计")
$newSheet   = $wb.Worksheets.Add($zongjiRef)
$newSheet.Name = "2022-Q1"

# NOTE: Worksheets.Add(Before) binds by position, so any handle obtained
# for "总计" *before* the insert now resolves to the freshly inserted
# sheet instead (indices shifted). Re-resolve "总计" by name afterwards.
$zongji = $wb.Worksheets.Item("总计")

# Copy header/data formatting from 2021-Q2 (3 populated rows: 1 header + 2 data)
$src.Range("A2:A3").Copy()
$newSheet.Range("A2:A3").PasteSpecial(-4122)  # xlPasteFormats
$src.Range("B1:H3").Copy()
$newSheet.Range("B1:H3").PasteSpecial(-4122)  # xlPasteFormats

# Header row
$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

# Row 2 - 004316
$newSheet.Range("A2").Value = 0
$newSheet.Range("B2").NumberFormat = "@"
$newSheet.Range("B2").Value = "004316"
$newSheet.Range("C2").Value = "前海开源沪港深裕鑫灵活配置混合A"
$newSheet.Range("D2").NumberFormat = "@"
$newSheet.Range("D2").Value = "0.64"
$newSheet.Range("E2").NumberFormat = "@"
$newSheet.Range("E2").Value = "90.55"
$newSheet.Range("F2").NumberFormat = "@"
$newSheet.Range("F2").Value = "3.08"
$newSheet.Range("G2").NumberFormat = "@"
$newSheet.Range("G2").Value = "0.0197"
$newSheet.Range("H2").Value = 7

# Row 3 - 004317
$newSheet.Range("A3").Value = 1
$newSheet.Range("B3").NumberFormat = "@"
$newSheet.Range("B3").Value = "004317"
$newSheet.Range("C3").Value = "前海开源沪港深裕鑫灵活配置混合C"
$newSheet.Range("D3").NumberFormat = "@"
$newSheet.Range("D3").Value = "0.47"
$newSheet.Range("E3").NumberFormat = "@"
$newSheet.Range("E3").Value = "90.55"
$newSheet.Range("F3").NumberFormat = "@"
$newSheet.Range("F3").Value = "3.08"
$newSheet.Range("G3").NumberFormat = "@"
$newSheet.Range("G3").Value = "0.0145"
$newSheet.Range("H3").Value = 7

# ------------------------------------------------------------------
# Step 2: add a new top data row to "总计" summarizing the 2022-Q1 quarter,
# pushing the existing 2021-Q2 / 2020-Q4 rows down by one.
# ------------------------------------------------------------------
$zongji.Rows.Item(2).Insert()

# Match formatting of the rows directly below/above the new blank row
$zongji.Range("A3").Copy()
$zongji.Range("A2").PasteSpecial(-4122)  # xlPasteFormats
$zongji.Range("B3:D3").Copy()
$zongji.Range("B2:D2").PasteSpecial(-4122)  # xlPasteFormats

$zongji.Range("A2").Value = 0
$zongji.Range("B2").Value = "2022-Q1"
$zongji.Range("C2").Value = 2
$zongji.Range("D2").Value = 0.03

$zongji.Range("A3").Value = 1
$zongji.Range("A4").Value = 2
